$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at position 305, shifting existing rows 305-380 down to 308-383
$ws.Range("A305:A307").EntireRow.Insert()

# Row 305 (new data)
$ws.Cells.Item(305, 1).Value = 5
$ws.Cells.Item(305, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(305, 3).Value = "Maule"
$ws.Cells.Item(305, 4).Value = 44551
$ws.Cells.Item(305, 5).Value = 7
$ws.Cells.Item(305, 6).Value = "Fruta"
$ws.Cells.Item(305, 7).Value = 100103
$ws.Cells.Item(305, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(305, 9).Value = 100103006
$ws.Cells.Item(305, 10).Value = "Nectarín"
$ws.Cells.Item(305, 11).Value = "Artic Star"
$ws.Cells.Item(305, 12).Value = "Primera"
$ws.Cells.Item(305, 13).Value = 500
$ws.Cells.Item(305, 14).Value = 12000
$ws.Cells.Item(305, 15).Value = 12000
$ws.Cells.Item(305, 16).Value = 12000
$ws.Cells.Item(305, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(305, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(305, 19).Value = 667
$ws.Cells.Item(305, 20).Value = 18

# Row 306 (new data)
$ws.Cells.Item(306, 1).Value = 5
$ws.Cells.Item(306, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(306, 3).Value = "Maule"
$ws.Cells.Item(306, 4).Value = 44551
$ws.Cells.Item(306, 5).Value = 7
$ws.Cells.Item(306, 6).Value = "Fruta"
$ws.Cells.Item(306, 7).Value = 100103
$ws.Cells.Item(306, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(306, 9).Value = 100103006
$ws.Cells.Item(306, 10).Value = "Nectarín"
$ws.Cells.Item(306, 11).Value = "Magique"
$ws.Cells.Item(306, 12).Value = "Especial"
$ws.Cells.Item(306, 13).Value = 200
$ws.Cells.Item(306, 14).Value = 15000
$ws.Cells.Item(306, 15).Value = 15000
$ws.Cells.Item(306, 16).Value = 15000
$ws.Cells.Item(306, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(306, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(306, 19).Value = 833
$ws.Cells.Item(306, 20).Value = 18

# Row 307 (new data)
$ws.Cells.Item(307, 1).Value = 5
$ws.Cells.Item(307, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(307, 3).Value = "Maule"
$ws.Cells.Item(307, 4).Value = 44551
$ws.Cells.Item(307, 5).Value = 7
$ws.Cells.Item(307, 6).Value = "Fruta"
$ws.Cells.Item(307, 7).Value = 100103
$ws.Cells.Item(307, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(307, 9).Value = 100103006
$ws.Cells.Item(307, 10).Value = "Nectarín"
$ws.Cells.Item(307, 11).Value = "Super Queen"
$ws.Cells.Item(307, 12).Value = "Especial"
$ws.Cells.Item(307, 13).Value = 500
$ws.Cells.Item(307, 14).Value = 14000
$ws.Cells.Item(307, 15).Value = 14000
$ws.Cells.Item(307, 16).Value = 14000
$ws.Cells.Item(307, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(307, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(307, 19).Value = 778
$ws.Cells.Item(307, 20).Value = 18
